# Update computed CG values on "GLOBAL RESULTS" and "LANDING GEARS" sheets
# (ACAnalysisManager re-run produced refreshed numbers).

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---
$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")

$wsGlobal.Range("C3").Value = 17.23823762216444
$wsGlobal.Range("C5").Value = -0.6999619377750967
$wsGlobal.Range("C7").Value = 34.04307493927112
$wsGlobal.Range("C9").Value = -17.95574380236645

$wsGlobal.Range("C13").Value = 16.40316226438204
$wsGlobal.Range("C15").Value = -0.7917816490286733
$wsGlobal.Range("C17").Value = 12.62133988368856
$wsGlobal.Range("C19").Value = -20.31114503534924

$wsGlobal.Range("C23").Value = 16.40316226438204
$wsGlobal.Range("C25").Value = -0.7917816490286733
$wsGlobal.Range("C27").Value = 12.62133988368856
$wsGlobal.Range("C29").Value = -20.31114503534924

$wsGlobal.Range("C33").Value = 16.40316226438204
$wsGlobal.Range("C35").Value = -0.7917816490286733
$wsGlobal.Range("C37").Value = 12.62133988368856
$wsGlobal.Range("C39").Value = -20.31114503534924

$wsGlobal.Range("C43").Value = 16.92417141603508
$wsGlobal.Range("C45").Value = -0.5851922612699278
$wsGlobal.Range("C47").Value = 25.986504960387446
$wsGlobal.Range("C49").Value = -15.011619563043268

$wsGlobal.Range("C53").Value = 16.840822886469397
$wsGlobal.Range("C55").Value = -0.6449790227816099
$wsGlobal.Range("C57").Value = 23.848410354925875
$wsGlobal.Range("C59").Value = -16.545296916828008

$wsGlobal.Range("C62").Value = 6.802712062372125
$wsGlobal.Range("C63").Value = 29.471094089476786

# --- LANDING GEARS sheet ---
$wsLanding = $wb.Worksheets.Item("LANDING GEARS")

$wsLanding.Range("C2").Value = 16.067908020990515
$wsLanding.Range("C4").Value = -3.129999999999999
$wsLanding.Range("C6").Value = 16.06790802099051
$wsLanding.Range("C8").Value = -4.979999999999999
$wsLanding.Range("C11").Value = 16.067908020990515
$wsLanding.Range("C14").Value = -3.129999999999999
